$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "post text" values for the two data rows.
$ws.Range("B2").Value = "Kuku 1"
$ws.Range("B3").Value = "Kuku 2"

# The old "Result"/"Pass" cell in D2 is removed entirely.
$ws.Range("D2").ClearContents()

# Move/selection moves to D4 after the edits.
$ws.Range("D4").Select()
